$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = 10
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 8
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = -5
$ws.Range("F13").Value = 2
